# Bumped hitter replacement level
#
# The original "replacement_hitters" sheet (plain numbers) is renamed to
# "unadjusted_replacement_hitters" and gets new "Adjustment" columns
# (H:M) holding per-position bumps to runs/hr/rbi/sb/avg.
#
# A brand new "replacement_hitters" sheet is inserted in front of it;
# it keeps the same Position/Runs/HR/RBI/SB/AVG layout but every stat
# cell is now a formula that adds the unadjusted value to the matching
# adjustment column on the other sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing data sheet and bolt on the adjustment columns.
# ---------------------------------------------------------------------
$unadj = $wb.Worksheets.Item("replacement_hitters")
$unadj.Name = "unadjusted_replacement_hitters"

$unadj.Range("H1").Value = "Adjustment"
$unadj.Range("I1").Value = "R"
$unadj.Range("J1").Value = "hr"
$unadj.Range("K1").Value = "rbi"
$unadj.Range("L1").Value = "sb"
$unadj.Range("M1").Value = "avg"

$adjustments = @{
    2  = @(7, 3, 7, 1.5, 0.02)
    3  = @(2, 1, 2, 0.5, 0.005)
    4  = @(3, 1.5, 3, 1, 0.01)
    5  = @(3, 1.5, 3, 1, 0.01)
    6  = @(2, 1, 2, 0.5, 0.005)
    7  = @(3, 1.5, 3, 1, 0.01)
    8  = @(2, 1, 2, 0.5, 0.005)
    9  = @(2, 1, 2, 0.5, 0.005)
    10 = @(2, 1, 2, 0.5, 0.005)
}

$adjCols = @("I", "J", "K", "L", "M")
foreach ($r in 2..10) {
    $vals = $adjustments[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $unadj.Cells.Item($r, 9 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 2. Insert the new "replacement_hitters" sheet in front of it.
# ---------------------------------------------------------------------
$adj = $wb.Worksheets.Add($unadj)
$adj.Name = "replacement_hitters"

$adj.Range("A1").Value = "Position"
$adj.Range("B1").Value = "Runs"
$adj.Range("C1").Value = "HR"
$adj.Range("D1").Value = "RBI"
$adj.Range("E1").Value = "SB"
$adj.Range("F1").Value = "AVG"

$positions = @{
    2  = "C"
    3  = "1b"
    4  = "2b"
    5  = "SS"
    6  = "3b"
    7  = "MI"
    8  = "CI"
    9  = "OF"
    10 = "Util"
}

$statCols = @("B", "C", "D", "E", "F")
foreach ($r in 2..10) {
    $adj.Cells.Item($r, 1).Value = $positions[$r]
    for ($i = 0; $i -lt $statCols.Length; $i++) {
        $c = $statCols[$i]
        $adjCol = $adjCols[$i]
        $adj.Range("$c$r").Formula = "=+unadjusted_replacement_hitters!$c$r+unadjusted_replacement_hitters!$adjCol$r"
    }
}

# ---------------------------------------------------------------------
# 3. Restore the view state: "unadjusted_replacement_hitters" is the
#    active/selected tab, with the selection left further down each
#    sheet (matching the author's last on-screen state).
# ---------------------------------------------------------------------
$adj.Range("D16").Select()

$unadjFresh = $wb.Worksheets.Item("unadjusted_replacement_hitters")
$unadjFresh.Range("H19").Select()
$unadjFresh.Activate()
